# Applies the "break out stock.yaml completed" edit:
#  1. day sheet: D627:D631 (bsecode) stored as inline-string text -> numeric
#  2. week sheet: append 13 new scraped rows (361:373), extends dimension to A1:I373
$wb = $excel.ActiveWorkbook

# --- 1. "day" sheet: fix bsecode column type (text -> number) for rows 627-631 ---
$dayWs = $wb.Worksheets.Item("day")
$dayWs.Cells.Item(627, 4).Value = 500820
$dayWs.Cells.Item(628, 4).Value = 532296
$dayWs.Cells.Item(629, 4).Value = 500302
$dayWs.Cells.Item(630, 4).Value = 539268
$dayWs.Cells.Item(631, 4).Value = 500253

# --- 2. "week" sheet: append new rows 361-373 ---
$weekWs = $wb.Worksheets.Item("week")

# bsecode (column D) must be kept as text in this newly-appended block -- force a
# text number format before assigning so numeric-looking codes like "500488" are
# not auto-converted to numbers.
$weekWs.Range("D361:D373").NumberFormat = "@"

$weekWs.Cells.Item(361, 1).Value = 1
$weekWs.Cells.Item(361, 2).Value = "ABBOTINDIA"
$weekWs.Cells.Item(361, 3).Value = "Abbott India Limited"
$weekWs.Cells.Item(361, 4).Value = "500488"
$weekWs.Cells.Item(361, 5).Value = -0.36
$weekWs.Cells.Item(361, 6).Value = 29165.5
$weekWs.Cells.Item(361, 7).Value = 5980
$weekWs.Cells.Item(361, 8).Value = "week"
$weekWs.Cells.Item(361, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(362, 1).Value = 2
$weekWs.Cells.Item(362, 2).Value = "LTTS"
$weekWs.Cells.Item(362, 3).Value = "L&t Technology Services Limited"
$weekWs.Cells.Item(362, 4).Value = "540115"
$weekWs.Cells.Item(362, 5).Value = -0.23
$weekWs.Cells.Item(362, 6).Value = 5344.25
$weekWs.Cells.Item(362, 7).Value = 98852
$weekWs.Cells.Item(362, 8).Value = "week"
$weekWs.Cells.Item(362, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(363, 1).Value = 3
$weekWs.Cells.Item(363, 2).Value = "TCS"
$weekWs.Cells.Item(363, 3).Value = "Tata Consultancy Services Limited"
$weekWs.Cells.Item(363, 4).Value = "532540"
$weekWs.Cells.Item(363, 5).Value = -0.93
$weekWs.Cells.Item(363, 6).Value = 4268.5
$weekWs.Cells.Item(363, 7).Value = 2503416
$weekWs.Cells.Item(363, 8).Value = "week"
$weekWs.Cells.Item(363, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(364, 1).Value = 4
$weekWs.Cells.Item(364, 2).Value = "CUMMINSIND"
$weekWs.Cells.Item(364, 3).Value = "Cummins India Limited"
$weekWs.Cells.Item(364, 4).Value = "500480"
$weekWs.Cells.Item(364, 5).Value = -1.52
$weekWs.Cells.Item(364, 6).Value = 3806.05
$weekWs.Cells.Item(364, 7).Value = 429408
$weekWs.Cells.Item(364, 8).Value = "week"
$weekWs.Cells.Item(364, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(365, 1).Value = 5
$weekWs.Cells.Item(365, 2).Value = "TORNTPHARM"
$weekWs.Cells.Item(365, 3).Value = "Torrent Pharmaceuticals Limited"
$weekWs.Cells.Item(365, 4).Value = "500420"
$weekWs.Cells.Item(365, 5).Value = -2.47
$weekWs.Cells.Item(365, 6).Value = 3396.55
$weekWs.Cells.Item(365, 7).Value = 335126
$weekWs.Cells.Item(365, 8).Value = "week"
$weekWs.Cells.Item(365, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(366, 1).Value = 6
$weekWs.Cells.Item(366, 2).Value = "LALPATHLAB"
$weekWs.Cells.Item(366, 3).Value = "Dr. Lal Path Labs Ltd."
$weekWs.Cells.Item(366, 4).Value = "539524"
$weekWs.Cells.Item(366, 5).Value = 0.94
$weekWs.Cells.Item(366, 6).Value = 3302.7
$weekWs.Cells.Item(366, 7).Value = 100909
$weekWs.Cells.Item(366, 8).Value = "week"
$weekWs.Cells.Item(366, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(367, 1).Value = 7
$weekWs.Cells.Item(367, 2).Value = "RELIANCE"
$weekWs.Cells.Item(367, 3).Value = "Reliance Industries Limited"
$weekWs.Cells.Item(367, 4).Value = "500325"
$weekWs.Cells.Item(367, 5).Value = -3.25
$weekWs.Cells.Item(367, 6).Value = 2953.15
$weekWs.Cells.Item(367, 7).Value = 13504407
$weekWs.Cells.Item(367, 8).Value = "week"
$weekWs.Cells.Item(367, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(368, 1).Value = 8
$weekWs.Cells.Item(368, 2).Value = "GRANULES"
$weekWs.Cells.Item(368, 3).Value = "Granules India Limited"
$weekWs.Cells.Item(368, 4).Value = "532482"
$weekWs.Cells.Item(368, 5).Value = 0.45
$weekWs.Cells.Item(368, 6).Value = 558.6
$weekWs.Cells.Item(368, 7).Value = 4010208
$weekWs.Cells.Item(368, 8).Value = "week"
$weekWs.Cells.Item(368, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(369, 1).Value = 9
$weekWs.Cells.Item(369, 2).Value = "COALINDIA"
$weekWs.Cells.Item(369, 3).Value = "Coal India Limited"
$weekWs.Cells.Item(369, 4).Value = "533278"
$weekWs.Cells.Item(369, 5).Value = -1.15
$weekWs.Cells.Item(369, 6).Value = 510.15
$weekWs.Cells.Item(369, 7).Value = 8881502
$weekWs.Cells.Item(369, 8).Value = "week"
$weekWs.Cells.Item(369, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(370, 1).Value = 10
$weekWs.Cells.Item(370, 2).Value = "INDUSTOWER"
$weekWs.Cells.Item(370, 3).Value = "Indus Towers Ltd (Bharti Infratel)"
$weekWs.Cells.Item(370, 4).Value = "534816"
$weekWs.Cells.Item(370, 5).Value = 0.04
$weekWs.Cells.Item(370, 6).Value = 392.55
$weekWs.Cells.Item(370, 7).Value = 9678928
$weekWs.Cells.Item(370, 8).Value = "week"
$weekWs.Cells.Item(370, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(371, 1).Value = 11
$weekWs.Cells.Item(371, 2).Value = "ASHOKLEY"
$weekWs.Cells.Item(371, 3).Value = "Ashok Leyland Limited"
$weekWs.Cells.Item(371, 4).Value = "500477"
$weekWs.Cells.Item(371, 5).Value = -1.73
$weekWs.Cells.Item(371, 6).Value = 235.4
$weekWs.Cells.Item(371, 7).Value = 10196690
$weekWs.Cells.Item(371, 8).Value = "week"
$weekWs.Cells.Item(371, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(372, 1).Value = 12
$weekWs.Cells.Item(372, 2).Value = "MANAPPURAM"
$weekWs.Cells.Item(372, 3).Value = "Manappuram Finance Limited"
$weekWs.Cells.Item(372, 4).Value = "531213"
$weekWs.Cells.Item(372, 5).Value = -0.64
$weekWs.Cells.Item(372, 6).Value = 201.35
$weekWs.Cells.Item(372, 7).Value = 2852007
$weekWs.Cells.Item(372, 8).Value = "week"
$weekWs.Cells.Item(372, 9).Value = "30/09/2024 18:34:57"

$weekWs.Cells.Item(373, 1).Value = 13
$weekWs.Cells.Item(373, 2).Value = "IDEA"
$weekWs.Cells.Item(373, 3).Value = "Idea Cellular Limited"
$weekWs.Cells.Item(373, 4).Value = "532822"
$weekWs.Cells.Item(373, 5).Value = -2.81
$weekWs.Cells.Item(373, 6).Value = 10.36
$weekWs.Cells.Item(373, 7).Value = 469780209
$weekWs.Cells.Item(373, 8).Value = "week"
$weekWs.Cells.Item(373, 9).Value = "30/09/2024 18:34:57"

Write-Output "applied stock.yaml break-out edit"